$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 "Save", matching the formatting used by the
# other header cells (e.g. G1): bold font, centered alignment, thin box border
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add new data cell H2 = 1 (plain numeric, no special style, like the
# other data cells in row 2)
$ws.Range("H2").Value = 1
